$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.027754068374634
$ws.Range("B1").Value = 2.196247816085815
$ws.Range("C1").Value = 7.043343067169189
$ws.Range("D1").Value = 2.348808288574219
$ws.Range("E1").Value = 1.325287342071533
